$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# Insert a single blank row before the totals block (old row 164) so the
# totals rows (164-167) shift down to 165-168 and their relative formulas
# (135-E164, 315-G164, ROUNDUP(C164/30,0)) get updated automatically.
$ws.Rows.Item(164).Insert()

# Row 162 was a blank spacer row; turn it into a new data row for the
# "Tablet neues Rezept anlegen screens" task.
$ws.Cells.Item(162, 1).Value = 22
$ws.Cells.Item(162, 2).Value = "Interface Design"
$ws.Cells.Item(162, 3).Value = "MockUps"
$ws.Cells.Item(162, 4).Value = "[FEATURE]"
$ws.Cells.Item(162, 5).Value = "Tablet neues Rezept anlegen screens"
$ws.Cells.Item(162, 6).Value = 44514
$ws.Cells.Item(162, 6).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(162, 7).Value = 44481
$ws.Cells.Item(162, 7).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(162, 9).Formula = "=ROUNDUP(((SUM(K162-J162)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(162, 9).NumberFormat = "0.00"
$ws.Cells.Item(162, 10).Value = 0.708333333333333
$ws.Cells.Item(162, 10).NumberFormat = "hh:mm"
$ws.Cells.Item(162, 11).Value = 0.791666666666667
$ws.Cells.Item(162, 11).NumberFormat = "hh:mm"

# Row 163 was also a blank spacer row; fill it the same way.
$ws.Cells.Item(163, 1).Value = 22
$ws.Cells.Item(163, 2).Value = "Interface Design"
$ws.Cells.Item(163, 3).Value = "MockUps"
$ws.Cells.Item(163, 4).Value = "[FEATURE]"
$ws.Cells.Item(163, 5).Value = "Tablet neues Rezept anlegen screens"
$ws.Cells.Item(163, 6).Value = 44515
$ws.Cells.Item(163, 6).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(163, 7).Value = 44481
$ws.Cells.Item(163, 7).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(163, 9).Formula = "=ROUNDUP(((SUM(K163-J163)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(163, 9).NumberFormat = "0.00"
$ws.Cells.Item(163, 10).Value = 0.541666666666667
$ws.Cells.Item(163, 10).NumberFormat = "hh:mm"
$ws.Cells.Item(163, 11).Value = 0.625
$ws.Cells.Item(163, 11).NumberFormat = "hh:mm"

# Clear the two stray empty L cells (L11/L12) that no longer carry content.
$ws.Cells.Item(11, 12).ClearContents()
$ws.Cells.Item(12, 12).ClearContents()

# The active selection moved to G164 after the edit.
$ws.Range("G164").Select()
